# Updated cryptos list (GitHub Actions refresh).
# Writes the new Price (col D) / Volume(1h) (col E) figures, and for rows
# 43-44 the Coin/Link also swap (OKB and dogwifhat trade places).
#
# Price/Volume cells are stored as TEXT in the workbook (e.g. "583.34",
# "  -0.24%  "), not numbers. Excel's Range.Value setter auto-detects
# numeric-looking strings and silently coerces them to real numbers
# (dropping formatting like trailing zeros, e.g. "583.34" -> 583.3400000...).
# To keep those cells textual we force NumberFormat="@" (Text) before the
# assignment, then reset the style back to "Normal" afterwards so no stray
# style/number-format sticks to the cell (matches the original, unstyled
# cells). Cells whose new value cannot parse as a plain number (e.g.
# "68.240.02", the percentages with their padding spaces, coin names,
# URLs) don't need this treatment - Excel already keeps those as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Force the cell to stay text even though $value parses as a number,
    # then strip the temporary Text format back off so the cell's style
    # is unchanged from before (no numFmt/style left behind).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = '68.240.02'
$ws.Range("E2").Value = '  -0.95%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '3.278.82'
$ws.Range("E3").Value = '  +0.33%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  -0.01%  '

# Row 5 - BNB
Set-TextValue "D5" '583.34'
$ws.Range("E5").Value = '  -0.24%  '

# Row 6 - Solana
Set-TextValue "D6" '185.14'
$ws.Range("E6").Value = '  +1.80%  '

# Row 7 - USDC
$ws.Range("E7").Value = '  +0.03%  '

# Row 8 - XRP
$ws.Range("E8").Value = '  +1.23%  '

# Row 9 - Dogecoin
Set-TextValue "D9" '0.131'
$ws.Range("E9").Value = '  -1.90%  '

# Row 10 - Toncoin
$ws.Range("E10").Value = '  -0.93%  '

# Row 11 - Cardano
$ws.Range("E11").Value = '  -2.94%  '

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '3.851.20'
$ws.Range("E12").Value = '  +0.42%  '

# Row 13 - TRON
Set-TextValue "D13" '0.139'
$ws.Range("E13").Value = '  +0.65%  '

# Row 14 - Avalanche
Set-TextValue "D14" '27.57'
$ws.Range("E14").Value = '  -3.18%  '

# Row 15 - WrappedBTC
$ws.Range("D15").Value = '68.228.71'
$ws.Range("E15").Value = '  -0.83%  '

# Row 16 - ShibaInu
$ws.Range("E16").Value = '  -1.37%  '

# Row 17 - WrappedEther
$ws.Range("D17").Value = '3.267.23'
$ws.Range("E17").Value = '  +1.46%  '

# Row 18 - Polkadot
Set-TextValue "D18" '5.73'
$ws.Range("E18").Value = '  -1.78%  '

# Row 19 - Chainlink
Set-TextValue "D19" '13.48'
$ws.Range("E19").Value = '  -0.43%  '

# Row 20 - BitcoinCash
Set-TextValue "D20" '417.88'
$ws.Range("E20").Value = '  +6.14%  '

# Row 21 - Uniswap
Set-TextValue "D21" '7.59'
$ws.Range("E21").Value = '  -1.02%  '

# Row 22 - Dai
$ws.Range("E22").Value = '  +0.22%  '

# Row 23 - Litecoin
Set-TextValue "D23" '71.48'
$ws.Range("E23").Value = '  -0.50%  '

# Row 24 - Polygon
$ws.Range("E24").Value = '  -0.82%  '

# Row 25 - PEPE
$ws.Range("E25").Value = '  -1.29%  '

# Row 26 - Kaspa
$ws.Range("E26").Value = '  -1.30%  '

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" '9.53'
$ws.Range("E27").Value = '  -0.65%  '

# Row 28 - Binance-PegBSC-USD
Set-TextValue "D28" '1.01'
$ws.Range("E28").Value = '  +0.72%  '

# Row 29 - PancakeSwap
$ws.Range("E29").Value = '  -1.60%  '

# Row 30 - EthereumClassic
$ws.Range("E30").Value = '  -0.86%  '

# Row 31 - NEARProtocol
Set-TextValue "D31" '5.50'
$ws.Range("E31").Value = '  -3.44%  '

# Row 32 - Aptos
$ws.Range("E32").Value = '  -2.74%  '

# Row 33 - Fetch.AI
$ws.Range("E33").Value = '  -2.34%  '

# Row 34 - USDe
$ws.Range("E34").Value = '  +0.05%  '

# Row 35 - Monero
Set-TextValue "D35" '164.08'
$ws.Range("E35").Value = '  +0.12%  '

# Row 36 - ImmutableX
$ws.Range("E36").Value = '  -2.75%  '

# Row 37 - Stacks
$ws.Range("E37").Value = '  -2.26%  '

# Row 38 - EnergySwap
Set-TextValue "D38" '27.38'
$ws.Range("E38").Value = '  +3.88%  '

# Row 39 - Mantle
Set-TextValue "D39" '0.801'
$ws.Range("E39").Value = '  -3.18%  '

# Row 40 - Filecoin
Set-TextValue "D40" '4.50'
$ws.Range("E40").Value = '  -1.80%  '

# Row 41 - RenderToken
$ws.Range("E41").Value = '  -3.85%  '

# Row 42 - Maker
$ws.Range("D42").Value = '2.672.10'
$ws.Range("E42").Value = '  +2.53%  '

# Row 43 - was dogwifhat, now OKB (rows 43/44 swap coin identity)
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D43" '40.88'
$ws.Range("E43").Value = '  -1.13%  '

# Row 44 - was OKB, now dogwifhat
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D44" '2.45'
$ws.Range("E44").Value = '  -1.08%  '

# Row 45 - Hedera
$ws.Range("E45").Value = '  -1.31%  '

# Row 46 - Bittensor
Set-TextValue "D46" '337.52'

# Row 47 - InjectiveProtocol
Set-TextValue "D47" '24.72'
$ws.Range("E47").Value = '  +0.04%  '

# Row 48 - VeChain
Set-TextValue "D48" '0.0276'
$ws.Range("E48").Value = '  -2.67%  '

# Row 49 - Cosmos
$ws.Range("E49").Value = '  +0.01%  '

# Row 50 - Stellar
$ws.Range("E50").Value = '  -1.07%  '

# Row 51 - ONDO
Set-TextValue "D51" '0.976'
$ws.Range("E51").Value = '  -0.65%  '
